$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the Price column (D) as text so numeric-looking strings are not
# auto-converted to numbers (the source data stores prices as inline strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "245.15"
$ws.Range("D3").Value = "25.21"
$ws.Range("D4").Value = "5.020"
$ws.Range("D5").Value = "0.05617"
$ws.Range("D6").Value = "6.569"
$ws.Range("D7").Value = "3.009"
$ws.Range("D8").Value = "0.8126"
$ws.Range("D9").Value = "0.8377"
$ws.Range("D10").Value = "0.1339"
$ws.Range("D11").Value = "0.06949"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.02839"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09394"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001514"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "0.009686"
$ws.Range("E15").Value = "14OneONEBestin24h"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006095"
$ws.Range("E16").Value = "15TigerCashTCH"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.499"
$ws.Range("E17").Value = "16LEOLEO"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.091"
$ws.Range("E18").Value = "17BTSETokenBTSE"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3184"
$ws.Range("E19").Value = "18BitpandaEcosystemTokenBEST"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "0.03263"
$ws.Range("E20").Value = "19LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("D21").Value = "0.1292"
$ws.Range("D22").Value = "3.736"
$ws.Range("D23").Value = "0.04661"
$ws.Range("D25").Value = "0.001243"
$ws.Range("D26").Value = "0.004535"
$ws.Range("D27").Value = "0.00009689"
$ws.Range("E27").Value = "26NitroExNTX"
$ws.Range("D28").Value = "0.0001938"
$ws.Range("D40").Value = "0.03664"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006223"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1057"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002726"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "0.008170"
$ws.Range("D45").Value = "0.00005287"
$ws.Range("D47").Value = "0.2257"
$ws.Range("D48").Value = "0.002035"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("D50").Value = "0.0001998"